$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The post "「私はサッカーが好きです」" (row 697) was removed entirely.
# Delete the whole row 697; all rows below it shift up by one, which
# matches the rest of the diff (rows 698..887 becoming 697..886) and the
# updated dimension (A1:C886).
$ws.Rows.Item(697).Delete()
